# Append 32 new guild-war result rows (A94:J125) to Sheet1, matching the
# columns: 방어팀, 방어팀 펫, 방어팀 스순, 공격팀, 공격팀 펫, 공격팀 스순, 속공, 날짜, 상대 길드, 기준
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    ,("트루드, 엘리시아, 파이", "유", "파1트2엘1", "오공, 엘리시아, 겔리두스", "크리", "겔1오2겔2", "후", 260105, "모현", "공격")
    ,("카일, 카구라, 파이", "이린", "파1카1카2", "트루드, 제이브, 카구라", "이린", "구2트2제1", "후", 260105, "모현", "공격")
    ,("트루드, 스파이크, 아멜리아", "리첼", "아1아2트2", "프레이야, 밀리아, 바네사", "연지", "바1프2프1", "선", 260105, "모현", "공격")
    ,("플라튼, 트루드, 아멜리아", "루", "아1아2트2", "크리스, 로지, 아킬라", "맬패로", "아2크2로1", "후", 260105, "모현", "공격")
    ,("카일, 카구라, 파이", "리첼", "구2카2카1", "바네사, 연희, 로지", "연지", "바1연2바2", "후", 260105, "모현", "공격")
    ,("카일, 카구라, 파이", "이린", "카2파1카1", "트루드, 카구라, 스파이크", "루", "구2스2트2", "후", 260105, "모현", "공격")
    ,("카일, 카구라, 파이", "이린", "카1파2카2", "키리엘, 밀리아, 프레이야", "연지", "밀2프2프1", "후", 260105, "모현", "공격")
    ,("오공, 겔리두스, 엘리시아", "파이크", "오2겔1엘2", "오공, 겔리두스, 엘리시아", "파이크", "겔2오2엘2", "선", 260105, "모현", "공격")
    ,("트루드, 스파이크, 아멜리아", "파이크", "아1아2트2", "트루드, 스파이크, 아멜리아", "파이크", "아2스2트2", "후", 260105, "모현", "공격")
    ,("트루드, 스파이크, 아멜리아", "파이크", "아1트2스2", "밀리아, 프레이야, 멜키르", "연지", "밀2프2프1", "선", 260105, "모현", "공격")
    ,("플라튼, 트루드, 아멜리아", "루", "아1아2트2", "플라튼, 엘리시아, 오공", "루", "엘2플2오2", "선", 260105, "모현", "공격")
    ,("플라튼, 콜트, 트루드", "카람", "플2콜1콜2", "크리스, 아킬라, 로지", "맬패로", "아2크1로1", "후", 260105, "모현", "공격")
    ,("트루드, 콜트, 엘리시아", "델로", "트2콜1콜2", "오공, 스파이크, 겔리두스", "파이크", "오2겔1오1", "후", 260105, "모현", "공격")
    ,("트루드, 스파이크, 키리엘", "카람", "트2스2키2", "오공, 겔리두스, 엘리시아", "크리", "겔2오2엘1", "선", 260105, "모현", "공격")
    ,("카일, 카구라, 파이", "이린", "카2파1카1", "플라튼, 바네사, 키리엘", "루", "바1바2키2", "후", 260105, "모현", "공격")
    ,("카일, 카구라, 파이", "이린", "카1카2파2", "오공, 겔리두스, 엘리시아", "파이크", "엘2오2오1", "후", 260105, "모현", "공격")
    ,("트루드, 스파이크, 아멜리아", "파라곤", "아2트2스2", "카일, 파이, 카구라", "이린", "카1파1카2", "선", 260105, "모현", "공격")
    ,("카일, 카구라, 파이", "이린", "파2카1카2", "실베스타, 쥬리, 키리엘", "리첼", "쥬2실1실2", "후", 260105, "모현", "공격")
    ,("카일, 카구라, 파이", "이린", "카1파2카2", "카일, 카구라, 파이", "이린", "카1카2파2", "후", 260105, "모현", "공격")
    ,("트루드, 스파이크, 아멜리아", "파이크", "아1트2아2", "트루드, 스파이크, 아멜리아", "파이크", "아2아1트2", "후", 260105, "모현", "공격")
    ,("플라튼, 트루드, 아멜리아", "루", "플2아1트2", "크리스, 녹스, 아킬라", "파이크", "아2크1크2", "선", 260105, "모현", "공격")
    ,("트루드, 스파이크, 아멜리아", "파이크", "아2아1트2", "프레이야, 바네사, 키리엘", "연지", "바1프2키2", "선", 260105, "모현", "공격")
    ,("트루드, 아멜리아, 엘리시아", "이린", "아1아2트2", "프레이야, 바네사, 키리엘", "연지", "바1프2키2", "선", 260105, "모현", "공격")
    ,("카일, 카구라, 파이", "이린", "구2카1카2", "트루드, 오공, 겔리두스", "파이크", "오2겔2트2", "후", 260105, "모현", "공격")
    ,("트루드, 스파이크, 아멜리아", "파이크", "아2트2스2", "트루드, 스파이크, 아멜리아", "파이크", "아1트2스2", "후", 260105, "모현", "공격")
    ,("카일, 카구라, 파이", "이린", "카1파2카2", "오공, 엘리시아, 겔리두스", "카람", "엘1오2겔2", "후", 260105, "모현", "공격")
    ,("카일, 카구라, 파이", "이린", "카2카1파2", "프레이야, 멜키르, 키리엘", "연지", "멜1프2키2", "후", 260105, "모현", "공격")
    ,("카일, 카구라, 파이", "이린", "구2카1카2", "트루드, 오공, 겔리두스", "이린", "겔1오2오1", "후", 260105, "모현", "공격")
    ,("카일, 카구라, 파이", "이린", "구2카1카2", "카일, 카구라, 파이", "크리", "구2카1카2", "후", 260105, "모현", "공격")
    ,("밀리아, 프레이야, 멜키르", "크리", "밀2프2프1", "밀리아, 프레이야, 쥬리", "크리", "쥬2프2프1", "선", 260105, "모현", "공격")
    ,("카일, 카구라, 파이", "이린", "구2카1카2", "카일, 카구라, 파이", "델로", "구2카1카2", "후", 260105, "모현", "공격")
    ,("카일, 카구라, 파이", "이린", "구2카1카2", "오공, 스파이크, 아멜리아", "맬패로", "아1오2스2", "후", 260105, "모현", "공격")
)

$startRow = 94
$r = $startRow
foreach ($row in $newRows) {
    for ($c = 1; $c -le 10; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $row[$c - 1]
        if ($c -eq 8) {
            $cell.Value = [double]$val
        } else {
            $cell.Value = [string]$val
        }
        $cell.HorizontalAlignment = -4108
    }
    $r++
}

$ws.Range("N120").Select()